$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1006.8410676762168
$ws.Range("C2").Value = 784.68221713555477
$ws.Range("D2").Value = 1103.9066269064899
$ws.Range("E2").Value = 1405.9551126350516
$ws.Range("F2").Value = 719.69678399746556
$ws.Range("G2").Value = 774.091548149593
$ws.Range("H2").Value = 1060.8049703138424
$ws.Range("I2").Value = 782.3268219189373
$ws.Range("J2").Value = 1130.9869802580927
$ws.Range("K2").Value = 1087.8628705064932
$ws.Range("L2").Value = 812.85936731041056
$ws.Range("M2").Value = 729.1864426670918
$ws.Range("N2").Value = 826.86839287623525
$ws.Range("O2").Value = 461.50722727523902
$ws.Range("P2").Value = 946.63924628826442
$ws.Range("Q2").Value = 500.3040003887192
$ws.Range("R2").Value = 689.23855025923604
$ws.Range("S2").Value = 1326.7332690296005
$ws.Range("T2").Value = 707.12147198043078
$ws.Range("U2").Value = 536.64590326082816
$ws.Range("V2").Value = 1089.7204951872588
$ws.Range("W2").Value = 967.72087494136485
$ws.Range("X2").Value = 763.76377689060121
$ws.Range("Y2").Value = 675.87984635763269
$ws.Range("Z2").Value = 273.31936521038307
$ws.Range("AA2").Value = 1658.9956138070947
$ws.Range("AB2").Value = 832.08217165400276
$ws.Range("AC2").Value = 898.8952163523918
$ws.Range("AD2").Value = 1005.3080529618902
$ws.Range("AE2").Value = 936.19748173051858
$ws.Range("AF2").Value = 838.91208182196033
$ws.Range("AG2").Value = 1101.6113265484894
$ws.Range("AH2").Value = 643.88711090781464
$ws.Range("AI2").Value = 1046.5595798839925
$ws.Range("AJ2").Value = 1537.9101646697975
$ws.Range("AK2").Value = 807.09777373235102
$ws.Range("AL2").Value = 921.14517117396565
$ws.Range("AM2").Value = 975.62009144007641
$ws.Range("AN2").Value = 417.24887649741385
$ws.Range("AO2").Value = 851.11674891164694
$ws.Range("AP2").Value = 454.03016517499162
$ws.Range("AQ2").Value = 476.22032922205875
$ws.Range("AR2").Value = 915.04243640297602
$ws.Range("AS2").Value = 1332.349838330452
$ws.Range("AT2").Value = 350.48076051474112
$ws.Range("AU2").Value = 1043.8642184216058
$ws.Range("AV2").Value = 992.86195296394544
$ws.Range("AW2").Value = 438.15613485493986
$ws.Range("AX2").Value = 854.45842571533592
$ws.Range("AY2").Value = 453.74051708991692
$ws.Range("B3").Value = 1256.6242850946503
$ws.Range("C3").Value = 899.75103386024921
$ws.Range("D3").Value = 528.38985073785295
$ws.Range("E3").Value = 1496.1471187305449
$ws.Range("F3").Value = 456.285361020219
$ws.Range("G3").Value = 1050.9346980209978
$ws.Range("H3").Value = 1072.7813803979839
$ws.Range("I3").Value = 589.6497493339499
$ws.Range("J3").Value = 1006.9601052355647
$ws.Range("K3").Value = 829.63542796963793
$ws.Range("L3").Value = 759.26891762922571
$ws.Range("M3").Value = 1036.1754179654108
$ws.Range("N3").Value = 770.17302558382835
$ws.Range("O3").Value = 444.87391105662761
$ws.Range("P3").Value = 1273.6451268954422
$ws.Range("Q3").Value = 585.28164379181419
$ws.Range("R3").Value = 799.87811287663681
$ws.Range("S3").Value = 795.74921583659318
$ws.Range("T3").Value = 738.77010299940105
$ws.Range("U3").Value = 545.17480982584766
$ws.Range("V3").Value = 756.4346694811212
$ws.Range("W3").Value = 542.46237332559338
$ws.Range("X3").Value = 517.10275689612502
$ws.Range("Y3").Value = 584.78716959042004
$ws.Range("Z3").Value = 346.35109702680694
$ws.Range("AA3").Value = 1426.1192511342895
$ws.Range("AB3").Value = 1108.1039680961651
$ws.Range("AC3").Value = 833.26334439337313
$ws.Range("AD3").Value = 1267.515929032468
$ws.Range("AE3").Value = 787.10236485966232
$ws.Range("AF3").Value = 1106.5721686571144
$ws.Range("AG3").Value = 1311.5007601466584
$ws.Range("AH3").Value = 1069.2692834284735
$ws.Range("AI3").Value = 1133.9687359503869
$ws.Range("AJ3").Value = 1626.883389688217
$ws.Range("AK3").Value = 1179.9377941935156
$ws.Range("AL3").Value = 974.68436770478218
$ws.Range("AM3").Value = 1377.297722886696
$ws.Range("AN3").Value = 641.67416989605465
$ws.Range("AO3").Value = 1653.9833540352877
$ws.Range("AP3").Value = 517.77112533875527
$ws.Range("AQ3").Value = 728.2778447998013
$ws.Range("AR3").Value = 689.50721370023564
$ws.Range("AS3").Value = 1130.6579010541188
$ws.Range("AT3").Value = 398.13473951382502
$ws.Range("AU3").Value = 1443.2658104534958
$ws.Range("AV3").Value = 892.0366397182255
$ws.Range("AW3").Value = 977.290564333936
$ws.Range("AX3").Value = 927.39974687635242
$ws.Range("AY3").Value = 634.35164670025995
